$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "Tabela1" (F2:I19) is re-sorted by the "Classificação" column
# (column G), descending, replacing the previous sort (which was by column I).
$lo = $ws.ListObjects.Item(1)
$sort = $lo.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("G2:G19"), 0, 2)
$sort.Header = 1
$sort.Apply()

# Cursor ends up on J7 after the sort/edit.
$ws.Range("J7").Select()
